$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Replace every "-" in column Q (Note column, rows 2-156) with the new
# placeholder note text "รออัพเดทจากธนาคาร".
for ($r = 2; $r -le 156; $r++) {
    $cell = $ws.Cells.Item($r, 17)
    if ($cell.Value2 -eq "-") {
        $cell.Value = "รออัพเดทจากธนาคาร"
    }
}

# Update the saved sheet view state (scroll position + active selection).
$ws.Application.ActiveWindow.ScrollRow = 58
$ws.Application.ActiveWindow.ScrollColumn = 11
$ws.Range("Q166").Select()
